$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link) - plain text, no numeric coercion needed ---
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

# --- Numeric-looking text columns (Price / Volume) ---
# Force text format so values like "308.38" / "-0.33%" are stored as strings,
# matching the source data (not auto-converted to Number/Percentage).
$deRanges = @('D2','E2','D3','E3','D4','E4','D5','E5','D6','E6','D7','E7','D8','E8','D9','E9','D10','E10','D11','E11','D12','E12','D13','E13','D14','E14','D15','E15','D16','E16','D17','E17','D19','E20','D21','E21','E22','D23','E23','D24','E24','D25','E25','D26','E27','D39','E39','D40','E40','D41','E41','D42','E42','D43','E43','D44','D45','E45','D46','E46','E47','E48','D49','E50','E51')
foreach ($r in $deRanges) { $ws.Range($r).NumberFormat = "@" }

$ws.Range('D2').Value = '308.38'
$ws.Range('E2').Value = '-0.33%'
$ws.Range('D3').Value = '39.36'
$ws.Range('E3').Value = '0.43%'
$ws.Range('D4').Value = '5.133'
$ws.Range('E4').Value = '0.43%'
$ws.Range('D5').Value = '0.08117'
$ws.Range('E5').Value = '-0.29%'
$ws.Range('D6').Value = '1.945'
$ws.Range('E6').Value = '-1.44%'
$ws.Range('D7').Value = '8.138'
$ws.Range('E7').Value = '2.64%'
$ws.Range('D8').Value = '0.9266'
$ws.Range('E8').Value = '-0.70%'
$ws.Range('D9').Value = '0.1416'
$ws.Range('E9').Value = '-0.13%'
$ws.Range('D10').Value = '0.1932'
$ws.Range('E10').Value = '-0.63%'
$ws.Range('D11').Value = '0.09089'
$ws.Range('E11').Value = '-0.54%'
$ws.Range('D12').Value = '0.03515'
$ws.Range('E12').Value = '0.57%'
$ws.Range('D13').Value = '0.09813'
$ws.Range('E13').Value = '-0.31%'
$ws.Range('D14').Value = '0.001397'
$ws.Range('E14').Value = '-0.82%'
$ws.Range('D15').Value = '0.005833'
$ws.Range('E15').Value = '-0.39%'
$ws.Range('D16').Value = '3.909'
$ws.Range('E16').Value = '9.20%'
$ws.Range('D17').Value = '4.235'
$ws.Range('E17').Value = '0.85%'
$ws.Range('D19').Value = '0.3453'
$ws.Range('E20').Value = '-0.89%'
$ws.Range('D21').Value = '4.791'
$ws.Range('E21').Value = '-0.35%'
$ws.Range('E22').Value = '-1.83%'
$ws.Range('D23').Value = '0.04380'
$ws.Range('E23').Value = '-1.80%'
$ws.Range('D24').Value = '0.001233'
$ws.Range('E24').Value = '-0.92%'
$ws.Range('D25').Value = '0.004838'
$ws.Range('E25').Value = '-0.41%'
$ws.Range('D26').Value = '0.0001303'
$ws.Range('E27').Value = '-9.93%'
$ws.Range('D39').Value = '0.02073'
$ws.Range('E39').Value = '-1.95%'
$ws.Range('D40').Value = '0.05102'
$ws.Range('E40').Value = '-0.75%'
$ws.Range('D41').Value = '0.007443'
$ws.Range('E41').Value = '-0.44%'
$ws.Range('D42').Value = '0.009836'
$ws.Range('E42').Value = '-1.96%'
$ws.Range('D43').Value = '0.1363'
$ws.Range('E43').Value = '-0.19%'
$ws.Range('D44').Value = '0.002134'
$ws.Range('D45').Value = '0.008489'
$ws.Range('E45').Value = '-16.21%'
$ws.Range('D46').Value = '0.00006399'
$ws.Range('E46').Value = '2.83%'
$ws.Range('E47').Value = '-0.26%'
$ws.Range('E48').Value = '-19.05%'
$ws.Range('D49').Value = '0.002560'
$ws.Range('E50').Value = '-0.26%'
$ws.Range('E51').Value = '-0.26%'

# Restore default (unstyled) cell style now that the text value is committed,
# so the saved style index matches the original (unstyled) cells.
foreach ($r in $deRanges) { $ws.Range($r).Style = "Normal" }
